$d = $word.ActiveDocument

function Get-ParagraphIndexForText($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return -1
    }
    $hitStart = $rng.Start
    $idx = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($hitStart -ge $p.Range.Start -and $hitStart -lt $p.Range.End) {
            $idx = $i
            break
        }
    }
    return $idx
}

# ------------------------------------------------------------------
# 1) Insert three new bullet paragraphs immediately before the bullet
#    "Developed and deployed custom analytical tools and algorithms..."
# ------------------------------------------------------------------
$newBullets = @(
    "• Developed meta-analytical techniques that identified systematic data quality issues across 20+ years of voter registration data",
    "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
    "• Created fraud detection systems analyzing 5+ terabyte datasets, uncovering demographic miscoding patterns across 2,000+ precincts"
)

$anchorIndex = Get-ParagraphIndexForText("Developed and deployed custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering")

$i = $anchorIndex
foreach ($t in $newBullets) {
    $p = $d.Paragraphs.Item($i)
    $ip = $p.Range.Duplicate
    $ip.Collapse(1)
    $ip.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($i)
    $newPara.Range.Text = $t
    $i = $i + 1
}

# ------------------------------------------------------------------
# 2) Remove the bullet
#    "Created fraud detection systems for campaign finance data analysis
#     across multi-terabyte datasets"
# ------------------------------------------------------------------
$removeIndex = Get-ParagraphIndexForText("Created fraud detection systems for campaign finance data analysis across multi-terabyte datasets")
if ($removeIndex -ne -1) {
    $d.Paragraphs.Item($removeIndex).Range.Delete()
}
